# Fix typo due to duplicates in venues
# Two duplicate rows need to be removed from the paper/venue table:
#   - Paper ID 2352968 (row 154) - duplicate "CAiSE Workshops" / "Workshop" entry
#   - Paper ID 2353177 (row 156) - duplicate "WETICE" / "Workshop" entry
# Removing the higher-numbered row first keeps the lower row index stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(156).Delete() | Out-Null
$ws.Rows.Item(154).Delete() | Out-Null

# Reflect the final cursor/selection position left after the cleanup
$ws.Range("A155:XFD155").Select() | Out-Null
